$p = $ppt.ActivePresentation

$newStyleId = "{6DACF13F-7D56-5426-5DBF-6772FC64968B}"

# --- Slide 1: table style id ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).Table.ApplyStyle($newStyleId)

# --- Slide 2: table style id + "y" run color ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).Table.ApplyStyle($newStyleId)

# Shape 11 on slide 2 is the textbox containing the single-letter run "y".
$yShape = $s2.Shapes.Item(11)
$yFont = $yShape.TextFrame.TextRange.Font
$yFont.Color.ObjectThemeColor = 4  # msoThemeColorLight2 (bg2)

# --- Slide 3: table style id ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).Table.ApplyStyle($newStyleId)
